$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("mapping_istanbul")
$ws2 = $wb.Worksheets.Item("mapping_Turkey")

# Update the industrial ("Sanayi") mapping row to switch CR/LFINF(...) and
# CR+PC/LWAL entries over to CR/LFM and CR+PC/LFM respectively.
# Order matters so the shared-string table is rebuilt in the same sequence
# as the authoritative edit.
$ws1.Range("M2").Value = "80% CR+PC/LFM+CDL/HBET:1-3`n20% CR+PC/LFM+CDL/HBET:4-6"
$ws1.Range("N2").Value = "100% CR+PC/LFM+CDL/HBET:1-3"
$ws1.Range("O2").Value = "100% CR+PC/LFM+CDL/HBET:1-3"
$ws1.Range("C2").Value = "60% CR/LFM+CDL/HBET:1-3`n40% CR/LFM+CDL/HBET:4-6"
$ws1.Range("E2").Value = "80% CR/LFM+CDL/HBET:1-3`n20% CR/LFM+CDL/HBET:4-6"
$ws1.Range("D2").Value = "100% CR/LFM+CDL/HBET:1-3"

$ws2.Range("C2").Value = "60% CR/LFM+CDL/HBET:1-3`n40% CR/LFM+CDL/HBET:4-6"
$ws2.Range("D2").Value = "100% CR/LFM+CDL/HBET:1-3"
$ws2.Range("E2").Value = "70% CR/LFM+CDL/HBET:1-3`n30% CR/LFM+CDL/HBET:4-6"
$ws2.Range("M2").Value = "100% CR+PC/LFM+CDL/HBET:1-3"
$ws2.Range("N2").Value = "100% CR+PC/LFM+CDL/HBET:1-3"
$ws2.Range("O2").Value = "100% CR+PC/LFM+CDL/HBET:1-3"

# Reflect the author's final cursor position/scroll state on each sheet.
$ws2.Activate()
$ws2.Range("B5").Select()

$ws1.Activate()
$ws1.Range("M2").Select()
